# Apply the commit: "add fr jurisdiction, snomed parameters, ..."
# - Set the Jurisdiction value (row 11, column B) on the Metadata sheet to "FRANCE"
# - Bump the Date value (row 8, column B) on the Metadata sheet to the new timestamp

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")

$metadata.Range("B8").Value = "2025-07-11T12:29:53+00:00"
$metadata.Range("B11").Value = "FRANCE"
